$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A-O_Demand")

# Rows 74,75,76,78,79,80,84,86,89 already have a "Gvkm" value in column AP (Target.Unit).
# All other data rows (2-90) need a new "nan" value written into column AP.
$gvkmRows = @(74,75,76,78,79,80,84,86,89)

for ($r = 2; $r -le 90; $r++) {
    if ($gvkmRows -contains $r) {
        $ws.Range("AP$r").Value = "Gvkm"
    } else {
        $ws.Range("AP$r").Value = "nan"
    }
}
